$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.089.57"
$ws.Range("E2").Value = "  -0.18%  "

$ws.Range("D3").Value = "1.877.52"
$ws.Range("E3").Value = "  -2.18%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "319.70"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.11%  "

$ws.Range("E6").Value = "  +0.25%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5036"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -3.50%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3959"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.23%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08212"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -4.00%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "42.05"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.33%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.093"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.12%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "23.57"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +5.15%  "

$ws.Range("D13").Value = "1.876.56"
$ws.Range("E13").Value = "  -2.31%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.295"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.89%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.192"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.01%  "

$ws.Range("E16").Value = "  +0.36%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "91.71"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.93%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001087"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.51%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06456"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.49%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.08"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.85%  "

$ws.Range("E21").Value = "  +0.24%  "

$ws.Range("D22").Value = "30.088.77"
$ws.Range("E22").Value = "  -0.22%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.836"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.94%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.15"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.25%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.154"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.47%  "

$ws.Range("D26").Value = "2.086.66"
$ws.Range("E26").Value = "  -2.53%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "161.21"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.81%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "21.11"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.00%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.246"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -8.18%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "127.49"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.30%  "

$ws.Range("E31").Value = "  -1.16%  "

$ws.Range("E32").Value = "  -2.65%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.920"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.34%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.696"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.66%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.02426"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.79%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.264"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.55%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06348"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.97%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2129"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.56%  "

$ws.Range("E39").Value = "  -5.07%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.490"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.69%  "

$ws.Range("B41").Value = "TheSandbox"
$ws.Range("C41").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6289"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.13%  "

$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.212"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.26%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "11.28"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.10%  "

$ws.Range("E44").Value = "  +0.25%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.15"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.33%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5903"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.16%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.088"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.32%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.628"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.56%  "

$ws.Range("E49").Value = "  -3.34%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "122.00"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.09%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "77.37"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.16%  "
